$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "{q3/4}{uoz} gin | {q3/4}{uoz} green Chartreuse | {q3/4}{uoz} maraschino liqueur | {q3/4}{uoz} lime juice | twist of lime"
$ws.Range("D1").Value = "Shake vigorously with ice. Strain into a {gcocktail glass} and garnish with lime twist."

$ws.Range("C2").Value = "{q1.5}{uoz} dry gin | {q1/2}{uoz} dry vermouth | {q1/2}{uoz} green Chartreuse | {q1}{utsp} absinthe | {q1-2} {udashes} orange bitters"
$ws.Range("D2").Value = "Shake with ice for a long time. Strain into a chilled {gcocktail glass}."

$ws.Range("C3").Value = "{q1} {ulime}, cut into eighths | {q1}{utsp} sugar | {q2}{uoz} cachaça"
$ws.Range("D3").Value = "Muddle lime and sugar in a {glowball glass} until the lime is juiced. Fill to brim with crushed ice and add cachaça. Garnish with sugar cane."

$ws.Range("C9").Select()
